# Updated symbol list with latest coinranking.com price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Coin) and C (Link) are plain text already; D (Price) and E (Volume)
# look numeric/percentage so we force them to Text format before writing, then
# drop back to the Normal style so no stray number format lingers on the cell.
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '305.98'
Set-TextCell $ws.Range("E2") '1.82%'

# Row 3
Set-TextCell $ws.Range("D3") '36.17'
Set-TextCell $ws.Range("E3") '-5.18%'

# Row 4
Set-TextCell $ws.Range("D4") '5.042'
Set-TextCell $ws.Range("E4") '1.27%'

# Row 5
Set-TextCell $ws.Range("D5") '0.07809'
Set-TextCell $ws.Range("E5") '1.17%'

# Row 6
Set-TextCell $ws.Range("D6") '2.128'
Set-TextCell $ws.Range("E6") '-3.30%'

# Row 7
Set-TextCell $ws.Range("D7") '7.911'
Set-TextCell $ws.Range("E7") '-0.71%'

# Row 8
Set-TextCell $ws.Range("B8") 'MXToken'
Set-TextCell $ws.Range("C8") 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws.Range("D8") '0.9181'
Set-TextCell $ws.Range("E8") '0.18%'

# Row 9
Set-TextCell $ws.Range("B9") 'LiechtensteinCryptoassetsExchange'
Set-TextCell $ws.Range("C9") 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell $ws.Range("D9") '0.09585'
Set-TextCell $ws.Range("E9") '5.96%'

# Row 10
Set-TextCell $ws.Range("B10") 'WazirX'
Set-TextCell $ws.Range("C10") 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell $ws.Range("D10") '0.1871'
Set-TextCell $ws.Range("E10") '4.12%'

# Row 11
Set-TextCell $ws.Range("B11") 'MandalaExchangeToken'
Set-TextCell $ws.Range("C11") 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell $ws.Range("D11") '0.08626'
Set-TextCell $ws.Range("E11") '1.78%'

# Row 12
Set-TextCell $ws.Range("B12") 'BitrueCoin'
Set-TextCell $ws.Range("C12") 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell $ws.Range("D12") '0.03509'
Set-TextCell $ws.Range("E12") '-1.04%'

# Row 13
Set-TextCell $ws.Range("B13") 'BitMartToken'
Set-TextCell $ws.Range("C13") 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell $ws.Range("D13") '0.09928'
Set-TextCell $ws.Range("E13") '-0.06%'

# Row 14
Set-TextCell $ws.Range("B14") 'BitForexToken'
Set-TextCell $ws.Range("C14") 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell $ws.Range("D14") '0.001434'
Set-TextCell $ws.Range("E14") '-3.08%'

# Row 15
Set-TextCell $ws.Range("B15") 'TigerCash'
Set-TextCell $ws.Range("C15") 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell $ws.Range("D15") '0.005694'
Set-TextCell $ws.Range("E15") '0.01%'

# Row 16
Set-TextCell $ws.Range("B16") 'LEO'
Set-TextCell $ws.Range("C16") 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws.Range("D16") '3.456'
Set-TextCell $ws.Range("E16") '-0.56%'

# Row 17
Set-TextCell $ws.Range("B17") 'GateToken'
Set-TextCell $ws.Range("C17") 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell $ws.Range("D17") '4.078'
Set-TextCell $ws.Range("E17") '2.14%'

# Row 18
Set-TextCell $ws.Range("D18") '2.399'
Set-TextCell $ws.Range("E18") '7.90%'

# Row 19
Set-TextCell $ws.Range("D19") '0.3419'
Set-TextCell $ws.Range("E19") '-1.27%'

# Row 20
Set-TextCell $ws.Range("D20") '0.1349'
Set-TextCell $ws.Range("E20") '2.48%'

# Row 21
Set-TextCell $ws.Range("D21") '4.760'
Set-TextCell $ws.Range("E21") '4.79%'

# Row 22
Set-TextCell $ws.Range("D22") '0.2297'
Set-TextCell $ws.Range("E22") '2.88%'

# Row 23
Set-TextCell $ws.Range("D23") '0.04614'
Set-TextCell $ws.Range("E23") '-1.03%'

# Row 24
Set-TextCell $ws.Range("E24") '14.88%'

# Row 25
Set-TextCell $ws.Range("D25") '0.001232'
Set-TextCell $ws.Range("E25") '0.11%'

# Row 26
Set-TextCell $ws.Range("D26") '0.0001400'
Set-TextCell $ws.Range("E26") '7.56%'

# Row 27
Set-TextCell $ws.Range("E27") '-42.63%'

# Row 39
Set-TextCell $ws.Range("D39") '0.01815'
Set-TextCell $ws.Range("E39") '4.30%'

# Row 40
Set-TextCell $ws.Range("D40") '0.04753'
Set-TextCell $ws.Range("E40") '1.48%'

# Row 41
Set-TextCell $ws.Range("E41") '-4.36%'

# Row 42
Set-TextCell $ws.Range("D42") '0.1400'
Set-TextCell $ws.Range("E42") '0.81%'

# Row 43
Set-TextCell $ws.Range("D43") '0.007752'
Set-TextCell $ws.Range("E43") '0.85%'

# Row 44
Set-TextCell $ws.Range("D44") '0.002239'
Set-TextCell $ws.Range("E44") '-2.35%'

# Row 45
Set-TextCell $ws.Range("D45") '0.01030'
Set-TextCell $ws.Range("E45") '5.33%'

# Row 46
Set-TextCell $ws.Range("D46") '0.00006246'
Set-TextCell $ws.Range("E46") '3.28%'

# Row 47
Set-TextCell $ws.Range("D47") '0.00000000752'
Set-TextCell $ws.Range("E47") '0.18%'

# Row 48
Set-TextCell $ws.Range("D48") '0.0005817'
Set-TextCell $ws.Range("E48") '0.29%'

# Row 49
Set-TextCell $ws.Range("D49") '28.79'
Set-TextCell $ws.Range("E49") '229.55%'

# Row 50
Set-TextCell $ws.Range("D50") '0.002006'
Set-TextCell $ws.Range("E50") '-25.66%'

# Row 51
Set-TextCell $ws.Range("D51") '0.00002107'
Set-TextCell $ws.Range("E51") '0.18%'
